# Fill in the remaining "Preconditions", "Method Inputs" and "Expected Result"
# columns of the automobile unit test plan (rows 7-12), completing the test
# cases for the Automobile class (__init__, __str__ and
# calculate_fuel_requirements).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$initInputsValid    = 'make = "HONDA",                                                       model = "CRV",                                   kilometers_per_litre =   20.0'
$initInputsBadMake   = 'make = "   ",                                                       model = "CRV",                                   kilometers_per_litre =   20.0'
$initInputsBadModel  = 'make = "HONDA",                                                       model = "   ",                                   kilometers_per_litre =   20.0'
$initInputsBadKpl    = 'make = "HONDA",                                                       model = "CRV",                                   kilometers_per_litre =   "meters"'
$strExpected         = '"Make: HONDA \n Model: CRV\nThis automobile can drive 20.0 kilometers per litre."'

# Row 7: __init__ - Attribute set to input values.
$ws.Range("G7").Value = "The automobile instance is created successfully with the attributes correctly set."
$ws.Range("E7").Value = "None"

# Rows 8-10: __init__ - Exception cases (preconditions all "None").
$ws.Range("E8").Value = "None"
$ws.Range("E9").Value = "None"
$ws.Range("E10").Value = "None"

$ws.Range("G8").Value = "ValueError"
$ws.Range("G9").Value = "ValueError"
$ws.Range("G10").Value = "ValueError"

# Row 10 method inputs (non-numeric kilometers_per_litre), reused as the
# precondition for the __str__ and calculate_fuel_requirements cases.
$ws.Range("F10").Value = $initInputsBadKpl
$ws.Range("E11").Value = $initInputsBadKpl
$ws.Range("E12").Value = $initInputsBadKpl

$ws.Range("F7").Value = $initInputsValid
$ws.Range("F8").Value = $initInputsBadMake
$ws.Range("F9").Value = $initInputsBadModel

# Row 11: __str__ - Returns string formatted appropriately
$ws.Range("G11").Value = $strExpected

# Row 12: calculate_fuel_requirements - Returns correct calculated value.
$ws.Range("F12").Value = "distance = 100.0"

# Row 11 method inputs (none needed for __str__).
$ws.Range("F11").Value = "None"

# Row 12 expected result is a plain number.
$ws.Range("G12").Value = 5
$ws.Range("G12").HorizontalAlignment = -4131

# Keep the current selection near the newly completed cells.
$null = $ws.Range("F11").Select()
